$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 284
$ws.Range("C3").Value = 174413
$ws.Range("C4").Value = 164406
$ws.Range("C8").Value = 64.31999999999999
